$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "Förändrad" (column C) for every existing data row (2-497)
#    from 45175 -> 45177.
$ws.Range("C2:C497").Value = 45177

# 2. Append the new row 498 with the new logging notification.
$ws.Cells.Item(498, 1).Value = "A 41511-2023"

$ws.Cells.Item(498, 2).Value = 45175
$ws.Cells.Item(498, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(498, 3).Value = 45177
$ws.Cells.Item(498, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(498, 4).Value = "JÖNKÖPINGS LÄN"
$ws.Cells.Item(498, 5).Value = "NÄSSJÖ"

$ws.Cells.Item(498, 7).Value = 1
$ws.Cells.Item(498, 8).Value = 0
$ws.Cells.Item(498, 9).Value = 0
$ws.Cells.Item(498, 10).Value = 0
$ws.Cells.Item(498, 11).Value = 0
$ws.Cells.Item(498, 12).Value = 0
$ws.Cells.Item(498, 13).Value = 0
$ws.Cells.Item(498, 14).Value = 0
$ws.Cells.Item(498, 15).Value = 0
$ws.Cells.Item(498, 16).Value = 0
$ws.Cells.Item(498, 17).Value = 0

$ws.Cells.Item(498, 18).Value = ""
$ws.Cells.Item(498, 18).WrapText = $true

# 3. Row 497 gets an explicit row height once a row is appended after it.
$ws.Rows.Item(497).RowHeight = 15
